$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Locator Value" column (C) for the rows whose CSS selectors were
# rescoped to the currently-visible Marketo form (form[style]:not(display:none))
# or rewritten to target inputs by their `name` attribute.
$ws.Cells.Item(25, 3).Value2 = "form[style]:not([style*=""display: none""]) .mktoCheckboxList input[name='termsandConditions']"
$ws.Cells.Item(26, 3).Value2 = "form[style]:not([style*=""display: none""]) .mktoCheckboxList input[name='Double_Opt_in_Compliant__c']"
$ws.Cells.Item(27, 3).Value2 = "form[style]:not([style*=""display: none""]) select#Country"
$ws.Cells.Item(28, 3).Value2 = "form[style]:not([style*=""display: none""]) select#jobTitle2"
$ws.Cells.Item(29, 3).Value2 = "form[style]:not([style*=""display: none""]) select#State"
$ws.Cells.Item(30, 3).Value2 = "form[style]:not([style*=""display: none""]) div#ValidMsgCountry"
$ws.Cells.Item(31, 3).Value2 = "form[style]:not([style*=""display: none""]) div#ValidMsgjobTitle2"
$ws.Cells.Item(32, 3).Value2 = "form[style]:not([style*=""display: none""]) .mktoError #ValidMsgState"
$ws.Cells.Item(33, 3).Value2 = "form[style]:not([style*=""display: none""]) div#ValidMsgCompany"
$ws.Cells.Item(34, 3).Value2 = "form[style]:not([style*=""display: none""]) div#ValidMsgEmail"
$ws.Cells.Item(35, 3).Value2 = "form[style]:not([style*=""display: none""]) div#ValidMsgFirstName"
$ws.Cells.Item(36, 3).Value2 = "form[style]:not([style*=""display: none""]) div#ValidMsgLastName"
$ws.Cells.Item(37, 3).Value2 = "form[style]:not([style*=""display: none""]) div#ValidMsgPhone"
$ws.Cells.Item(38, 3).Value2 = "form[style]:not([style*=""display: none""]) .mktoError #ValidMsgPostalCode"
$ws.Cells.Item(40, 3).Value2 = "form[style]:not([style*=""display: none""]) input#Company"
$ws.Cells.Item(41, 3).Value2 = "form[style]:not([style*=""display: none""]) input#Email"
$ws.Cells.Item(42, 3).Value2 = "form[style]:not([style*=""display: none""]) input#FirstName"
$ws.Cells.Item(43, 3).Value2 = "form[style]:not([style*=""display: none""]) input#LastName"
$ws.Cells.Item(44, 3).Value2 = "form[style]:not([style*=""display: none""]) input#Phone"
$ws.Cells.Item(45, 3).Value2 = "form[style]:not([style*=""display: none""]) input#PostalCode"
$ws.Cells.Item(47, 3).Value2 = "form[style]:not([style*=""display: none""]) div.mktoError #ValidMsgtermsandConditions"

# Reflect the author's final cursor position/selection in the sheet view.
$ws.Range("C40").Select() | Out-Null
